$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing order to reflect new status ---
$ws.Range("C2").Value = "Yoff"
$ws.Range("D2").Value = "05h26"
$ws.Range("E2").Value = "Bazoungoula"
$ws.Range("F2").Value = "En attente"
$ws.Range("G2").Value = "George N'gock"
$ws.Range("H2").Value = "Chez le livreur"
$ws.Range("I2").Value = "Vosgienne Taille Petite, Vosgienne - Antoinette Taille Grande"
$ws.Range("J2").Value = 1200
$ws.Range("K2").Value = 7500
$ws.Range("L2").Value = 8700

# --- Row 3: replaced by a new "commandée sur place" order ---
$ws.Range("A3").Value = "27 décembre 2024"
$ws.Range("B3").Value = "Oui"
$ws.Range("C3").Value = "Commandée sur place"
$ws.Range("D3").Value = "Heure sur place: 12:09"
$ws.Range("E3").Value = "Commandée sur place"
$ws.Range("F3").Value = "Sur place"
$ws.Range("G3").Value = "Coulibaly Yelanto"
$ws.Range("H3").Value = "Chez Izoua"
$ws.Range("I3").Value = "Antoinette Taille Petite"
$ws.Range("J3").Value = "Commandée sur place"
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2500

# --- Row 4: new delivered order, previously-existing data moved here ---
$ws.Range("A4").Value = "27 décembre 2024"
$ws.Range("B4").Value = "Non"
$ws.Range("C4").Value = "Ouest Foire Dakar"
$ws.Range("D4").Value = "02h51"
$ws.Range("E4").Value = "Cagil"
$ws.Range("F4").Value = "Livré"
$ws.Range("G4").Value = "yve kate"
$ws.Range("H4").Value = "Chez le livreur"
$ws.Range("I4").Value = "Antoinette Taille Grande"
$ws.Range("J4").Value = 900
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 5900
